# Thomas Brinson 2026-01-19 timesheet export: full-week coverage + OT row,
# new employees on the daily log, persisted "log" sheet row, and a
# re-issued employee id (simulator re-roll).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Weekly Timesheet"
$ws2 = $wb.Worksheets.Item(2)   # "Jason Schema" (flattened log)

# ---------------------------------------------------------------------
# 1) Insert a new row under the existing 5 data rows (row 7) so the
#    SUBTOTAL / HOURLY SUBTOTAL / ADMIN SUBTOTAL / GRAND TOTAL block
#    shifts down by one row, exactly like Excel's own Insert does. The
#    new row inherits the plain (unstyled) look of the data rows above
#    it because columns E/F already carry the currency style.
# ---------------------------------------------------------------------
$ws1.Rows(7).Insert()

# ---------------------------------------------------------------------
# 2) Update the five original daily rows: new client names, and hours
#    bumped from 8 -> 9 (except the last day, split into Reg 4 / OT 5),
#    rate bumped from 0 -> 95 with totals recalculated accordingly.
# ---------------------------------------------------------------------
$ws1.Range("B2").Value = "McClure"
$ws1.Range("C2").Value = 9
$ws1.Range("E2").Value = 95
$ws1.Range("F2").Value = 855

$ws1.Range("B3").Value = "Evans"
$ws1.Range("C3").Value = 9
$ws1.Range("E3").Value = 95
$ws1.Range("F3").Value = 855

$ws1.Range("B4").Value = "Fritts"
$ws1.Range("C4").Value = 9
$ws1.Range("E4").Value = 95
$ws1.Range("F4").Value = 855

$ws1.Range("B5").Value = "Hendricks"
$ws1.Range("C5").Value = 9
$ws1.Range("E5").Value = 95
$ws1.Range("F5").Value = 855

$ws1.Range("B6").Value = "Regan"
$ws1.Range("C6").Value = 4
$ws1.Range("E6").Value = 95
$ws1.Range("F6").Value = 380

# ---------------------------------------------------------------------
# 3) Populate the newly-inserted row 7: same date/client as the last
#    day, but an OT line (5 hours @ 95 = 712.5).
# ---------------------------------------------------------------------
$ws1.Range("A7").Value = "2026-01-23"
$ws1.Range("B7").Value = "Regan"
$ws1.Range("C7").Value = 5
$ws1.Range("D7").Value = "OT"
$ws1.Range("E7").Value = 95
$ws1.Range("F7").Value = 712.5

# ---------------------------------------------------------------------
# 4) SUBTOTAL row (now row 9): total hours 45, Reg 40 / OT 5 summary,
#    and the grand hourly total of 4512.5.
# ---------------------------------------------------------------------
$ws1.Range("C9").Value = 45
$ws1.Range("D9").Value = "Reg: 40 / OT: 5"
$ws1.Range("F9").Value = 4512.5

# ---------------------------------------------------------------------
# 5) HOURLY SUBTOTAL row (now row 12) carries the same grand total.
#    ADMIN SUBTOTAL (row 13) / GRAND TOTAL (row 14) keep their existing
#    totals (0 and 4512.5 respectively) which fall out of the shift.
# ---------------------------------------------------------------------
$ws1.Range("F12").Value = 4512.5

# ---------------------------------------------------------------------
# 6) Mirror the same data into the flattened "Jason Schema" log sheet:
#    refresh the five existing rows and append the new OT row.
# ---------------------------------------------------------------------
$ws2.Range("D2").Value = "McClure"
$ws2.Range("E2").Value = 9
$ws2.Range("F2").Value = 95
$ws2.Range("G2").Value = 855

$ws2.Range("D3").Value = "Evans"
$ws2.Range("E3").Value = 9
$ws2.Range("F3").Value = 95
$ws2.Range("G3").Value = 855

$ws2.Range("D4").Value = "Fritts"
$ws2.Range("E4").Value = 9
$ws2.Range("F4").Value = 95
$ws2.Range("G4").Value = 855

$ws2.Range("D5").Value = "Hendricks"
$ws2.Range("E5").Value = 9
$ws2.Range("F5").Value = 95
$ws2.Range("G5").Value = 855

$ws2.Range("D6").Value = "Regan"
$ws2.Range("E6").Value = 4
$ws2.Range("F6").Value = 95
$ws2.Range("G6").Value = 380

$ws2.Range("A7").Value = "Thomas Brinson"
$ws2.Range("B7").Value = "emp_4nlnrvy7"
$ws2.Range("C7").Value = "2026-01-23"
$ws2.Range("D7").Value = "Regan"
$ws2.Range("E7").Value = 5
$ws2.Range("F7").Value = 95
$ws2.Range("G7").Value = 712.5
$ws2.Range("H7").Value = "OT"
$ws2.Range("I7").Value = ""

# ---------------------------------------------------------------------
# 7) Employee id was re-issued by the simulator; update everywhere it
#    appears (header sheet log rows B2:B6 above + B7 already set).
# ---------------------------------------------------------------------
$ws2.Range("B2").Value = "emp_4nlnrvy7"
$ws2.Range("B3").Value = "emp_4nlnrvy7"
$ws2.Range("B4").Value = "emp_4nlnrvy7"
$ws2.Range("B5").Value = "emp_4nlnrvy7"
$ws2.Range("B6").Value = "emp_4nlnrvy7"

Write-Host "edit applied"
